# Prepped for revamp of simulated stock
# - Rename "Coupon Yield"/"Sale at Maturity"/"Maturity" shared labels to the
#   new "Dividends"/"Capital Appreciation" scheme used by the revamped
#   simulated-stock model, and refresh the simulated Data sheet + ledgers to
#   match the new month-by-month projection.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsRevenue = $wb.Worksheets.Item("Revenue Ledger")
$wsInvestments = $wb.Worksheets.Item("Investments Ledger")
$wsCapitalGains = $wb.Worksheets.Item("Capital Gains Ledger")

# --- Data sheet: refreshed simulated-stock projection (rows 2-15) ---------
$wsData.Range("F2").Value = 1.0
$wsData.Range("I2").Value = 11.0
$wsData.Range("J2").Value = 11.0
$wsData.Range("K2").Value = 10.0
$wsData.Range("L2").Value = 1.0
$wsData.Range("M2").Value = 1.8409443766154858
$wsData.Range("F3").Value = 1.0
$wsData.Range("I3").Value = 12.0
$wsData.Range("J3").Value = 12.0
$wsData.Range("K3").Value = 10.0
$wsData.Range("L3").Value = 1.0
$wsData.Range("M3").Value = 1.6130352902246758
$wsData.Range("F4").Value = 1.0
$wsData.Range("I4").Value = 13.0
$wsData.Range("J4").Value = 13.0
$wsData.Range("K4").Value = 10.0
$wsData.Range("L4").Value = 1.0
$wsData.Range("M4").Value = 1.4334151101796073
$wsData.Range("F5").Value = 1.0
$wsData.Range("I5").Value = 14.0
$wsData.Range("J5").Value = 14.0
$wsData.Range("K5").Value = 10.0
$wsData.Range("L5").Value = 1.0
$wsData.Range("M5").Value = 1.2885409112665345
$wsData.Range("F6").Value = 1.0
$wsData.Range("I6").Value = 15.0
$wsData.Range("J6").Value = 15.0
$wsData.Range("K6").Value = 10.0
$wsData.Range("L6").Value = 1.0
$wsData.Range("M6").Value = 1.1694252129716127
$wsData.Range("F7").Value = 1.0
$wsData.Range("I7").Value = 16.0
$wsData.Range("J7").Value = 16.0
$wsData.Range("K7").Value = 10.0
$wsData.Range("L7").Value = 1.0
$wsData.Range("M7").Value = 1.0698899917795224
$wsData.Range("F8").Value = 1.0
$wsData.Range("I8").Value = 17.0
$wsData.Range("J8").Value = 17.0
$wsData.Range("K8").Value = 10.0
$wsData.Range("L8").Value = 1.0
$wsData.Range("M8").Value = 0.9855599520654272
$wsData.Range("F9").Value = 1.0
$wsData.Range("I9").Value = 18.0
$wsData.Range("J9").Value = 18.0
$wsData.Range("K9").Value = 10.0
$wsData.Range("L9").Value = 1.0
$wsData.Range("M9").Value = 0.9132562918007303
$wsData.Range("F10").Value = 1.0
$wsData.Range("I10").Value = 19.0
$wsData.Range("J10").Value = 19.0
$wsData.Range("K10").Value = 10.0
$wsData.Range("L10").Value = 1.0
$wsData.Range("M10").Value = 0.8506178062217085
$wsData.Range("F11").Value = 1.0
$wsData.Range("I11").Value = 20.0
$wsData.Range("J11").Value = 20.0
$wsData.Range("K11").Value = 10.0
$wsData.Range("L11").Value = 1.0
$wsData.Range("M11").Value = 0.7958563260221301
$wsData.Range("F12").Value = 1.0
$wsData.Range("I12").Value = 21.0
$wsData.Range("J12").Value = 21.0
$wsData.Range("K12").Value = 10.0
$wsData.Range("L12").Value = 1.0
$wsData.Range("M12").Value = 0.7475943544285117
$wsData.Range("E13").Value = 0.0
$wsData.Range("F13").Value = 1.0
$wsData.Range("G13").Value = 0.0
$wsData.Range("I13").Value = 22.0
$wsData.Range("J13").Value = 22.0
$wsData.Range("K13").Value = 10.0
$wsData.Range("L13").Value = 1.0
$wsData.Range("M13").Value = 0.7047545660620107
$wsData.Range("F14").Value = 1.0
$wsData.Range("I14").Value = 23.0
$wsData.Range("J14").Value = 23.0
$wsData.Range("K14").Value = 10.0
$wsData.Range("L14").Value = 1.0
$wsData.Range("M14").Value = 0.6664829255615827
$wsData.Range("F15").Value = 1.0
$wsData.Range("I15").Value = 24.0
$wsData.Range("J15").Value = 24.0
$wsData.Range("K15").Value = 10.0
$wsData.Range("L15").Value = 1.0
$wsData.Range("M15").Value = 0.6320941327229255

# --- Revenue Ledger: rename "Coupon Yield" -> "Dividends" and add two more
#     "Dividends" entries (month columns 2 and 8) to match the new model ---
$wsRevenue.Range("G2").Value = "Dividends"
$wsRevenue.Range("H2").Value = 0.0
$wsRevenue.Range("P2").Value = "Dividends"
$wsRevenue.Range("Q2").Value = 0.0
$wsRevenue.Range("Y2").Value = "Dividends"
$wsRevenue.Range("Z2").Value = 0.0
$wsRevenue.Range("AH2").Value = "Dividends"
$wsRevenue.Range("AI2").Value = 0.0

# --- Investments Ledger: drop the old "Sale at Maturity" row entirely -----
$wsInvestments.Rows.Item(2).Delete()

# --- Capital Gains Ledger: replace the single "Maturity" entry with one
#     "Capital Appreciation" entry per simulated month -------------------
$capitalGainsPairs = @(
    @("A2", "B2"), @("D2", "E2"), @("G2", "H2"), @("J2", "K2"), @("M2", "N2"),
    @("P2", "Q2"), @("S2", "T2"), @("V2", "W2"), @("Y2", "Z2"), @("AB2", "AC2"),
    @("AE2", "AF2"), @("AH2", "AI2"), @("AK2", "AL2"), @("AN2", "AO2")
)
foreach ($pair in $capitalGainsPairs) {
    $wsCapitalGains.Range($pair[0]).Value = "Capital Appreciation"
    $wsCapitalGains.Range($pair[1]).Value = 1.0
}
